$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 21.85308466666666
$ws.Range("N2").Value = 65.559254
$ws.Range("O2").Value = 0.407053040353553
$ws.Range("P2").Value = 0.407053040353553
$ws.Range("Q2").Value = 4.097009028945111
$ws.Range("R2").Value = 36.873081260506
$ws.Range("S2").Value = 0.407053040353553
$ws.Range("T2").Value = 0.407053040353553

# Row 3
$ws.Range("O3").Value = 0.1342711086924142
$ws.Range("P3").Value = 0.1342711086924142
$ws.Range("S3").Value = 0.1342711086924142
$ws.Range("T3").Value = 0.1342711086924142

# Row 4
$ws.Range("M4").Value = 11.375406
$ws.Range("N4").Value = 34.126218
$ws.Range("O4").Value = 0.2118874139822907
$ws.Range("P4").Value = 0.2118874139822907
$ws.Range("Q4").Value = 2.132657325078
$ws.Range("R4").Value = 19.193915925702
$ws.Range("S4").Value = 0.2118874139822907
$ws.Range("T4").Value = 0.2118874139822907

# Row 5
$ws.Range("M5").Value = 3.401340666666667
$ws.Range("N5").Value = 10.204022
$ws.Range("O5").Value = 0.06335609277882483
$ws.Range("P5").Value = 0.06335609277882483
$ws.Range("Q5").Value = 0.6376822144064445
$ws.Range("R5").Value = 5.739139929658
$ws.Range("S5").Value = 0.06335609277882483
$ws.Range("T5").Value = 0.06335609277882483

# Row 6
$ws.Range("M6").Value = 9.847764666666666
$ws.Range("N6").Value = 29.543294
$ws.Range("O6").Value = 0.1834323441929172
$ws.Range("P6").Value = 0.1834323441929172
$ws.Range("Q6").Value = 1.846255637118445
$ws.Range("R6").Value = 16.616300734066
$ws.Range("S6").Value = 0.1834323441929172
$ws.Range("T6").Value = 0.1834323441929172
